$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a literal TEXT value into a cell without Excel's automatic
# number/date inference (which would otherwise turn digit-only or
# date-shaped strings into Number/Date cells and touch NumberFormat/style).
# We do this by writing a string-literal formula ("=""...""") and then
# collapsing it back down to a plain value via copy / paste-special-values,
# which keeps the cell's existing style (s="2") untouched.
function Set-TextValue([string]$addr, [string]$text) {
    $r = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

Set-TextValue "F2"  "5653094250"
Set-TextValue "N2"  "2024-05-28"
Set-TextValue "O2"  "02:35:55 PM"
Set-TextValue "P2"  "2024-05-31 05:00:00 PM"
Set-TextValue "AC2" "2024-05-28"
Set-TextValue "AE2" "6571600101"
Set-TextValue "AK2" "3"
Set-TextValue "AN2" "126178"
Set-TextValue "AT2" "3756905950"
Set-TextValue "AX2" "4043615867"

$excel.CutCopyMode = $false
